$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The writer used to append the whole record block again before the final
# summary row - so re-insert two more copies of rows 1-8 ahead of the
# trailing row (previously row 9, the "Britney" record), pushing it down
# to the end of the table (row 25).

$blockRows = 8
$copies = 2

# Make room: insert (copies * blockRows) blank rows right before the
# trailing row so it slides down to its new position.
for ($i = 0; $i -lt ($blockRows * $copies); $i++) {
    $ws.Rows.Item(9).Insert()
}

# Re-write the source block (rows 1-8) twice into the freshly inserted
# rows 9-16 and 17-24.
for ($c = 0; $c -lt $copies; $c++) {
    for ($r = 1; $r -le $blockRows; $r++) {
        $destRow = 8 + ($c * $blockRows) + $r

        $ws.Cells.Item($destRow, 1).Value = $ws.Cells.Item($r, 1).Value()
        $ws.Cells.Item($destRow, 2).Value = $ws.Cells.Item($r, 2).Value()
        $ws.Cells.Item($destRow, 3).Value = $ws.Cells.Item($r, 3).Value()
        $ws.Cells.Item($destRow, 4).Value = $ws.Cells.Item($r, 4).Value()
        $ws.Cells.Item($destRow, 5).Value = $ws.Cells.Item($r, 5).Value()
    }
}
